# Hana_T85.xlsx edit:
# - Remove the stray "pickup phone order" test-case row value in C2
#   (collapse it down to the same "Abish David" value/style as B2,
#   which drops the now-unused shared string + cell style).
# - Narrow column C now that it no longer holds the long search string.
# - Move the active selection from B3 to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 currently holds the long "Abish David 114 N CHURCH ST PICK UP (...)"
# search-and-select string; replace it with the plain customer name that's
# already in B2, same as the other "Dispatch" test rows.
$ws.Range("C2").Value2 = $ws.Range("B2").Value2

# B2 and C2 should carry the same cell style as A2 (border, no quote-prefix
# formatting) instead of the separate "applyBorder only" style, so the old
# style record drops out entirely.
$ws.Range("A2").Copy()
$ws.Range("B2:C2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column C no longer needs to fit the long search string - shrink it to
# match the shorter "Abish David" contents.
$ws.Columns("C").ColumnWidth = 17.15

# Selection moves from B3 to A3.
$ws.Range("A3").Select()
